# filter date v1 , removed unnecessary filters
#
# The "Generate Yearly Report for Vendor / WI4 / Open" rows are replaced
# with a new set of "Research Client Check Copy / WI2 / Open" rows, and the
# table grows from 12 data rows (rows 2-13) to 28 data rows (rows 2-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: WIID (col B) and Date-serial (col F) for each row, in order.
# Description (col C) = "Research Client Check Copy"
# Type        (col D) = "WI2"
# Status      (col E) = "Open"
$rows = @(
    @(647572, 42788),
    @(587242, 43063),
    @(962522, 43441),
    @(612812, 43075),
    @(837862, 43242),
    @(906912, 43381),
    @(490572, 43077),
    @(295372, 43051),
    @(515842, 43536),
    @(699802, 42888),
    @(251502, 43064),
    @(321742, 42823),
    @(239522, 43004),
    @(802782, 43150),
    @(636492, 43516),
    @(697552, 43231),
    @(760052, 42896),
    @(191382, 43141),
    @(479542, 43279),
    @(758992, 43181),
    @(746282, 43188),
    @(562152, 42867),
    @(436842, 42803),
    @(583912, 43442),
    @(367012, 43319),
    @(982202, 42774),
    @(419692, 43455),
    @(352612, 43085)
)

# Clear out the previous data block (rows 2-13) before re-populating.
$ws.Range("A2:F13").ClearContents()

$startRow = 2
$lastRow = $startRow + $rows.Count - 1

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = "Research Client Check Copy"
    $ws.Cells.Item($r, 4).Value = "WI2"
    $ws.Cells.Item($r, 5).Value = "Open"
    $ws.Cells.Item($r, 6).Value = $row[1]
    $r = $r + 1
}

# Column F carries a short-date display format (style index 1 in the
# original file) - propagate that format to the newly written date cells
# by copying it from the existing, already-formatted F2 cell.
$ws.Cells.Item(2, 6).Copy()
$ws.Range("F" + $startRow + ":F" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output ("Wrote " + $rows.Count + " rows (" + $startRow + ":" + $lastRow + ")")
